$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.369
$ws.Range("E4").Value = 16.11029999999999
$ws.Range("D7").Value = -6.903800000000003
$ws.Range("A8").Value = -22.43020000000002
$ws.Range("A10").Value = -21.85209999999999
$ws.Range("E11").Value = 16.5645
$ws.Range("A12").Value = -21.5429
$ws.Range("D14").Value = -8.111300000000002
$ws.Range("E14").Value = 16.44009999999998
$ws.Range("D15").Value = -7.906200000000002
$ws.Range("A18").Value = -22.41870000000001
$ws.Range("D18").Value = -8.149599999999994
$ws.Range("E18").Value = 16.5745
$ws.Range("E19").Value = 16.47450000000001
$ws.Range("D20").Value = -7.843399999999998
$ws.Range("E21").Value = 16.636
$ws.Range("A25").Value = -21.68639999999998
$ws.Range("E27").Value = 16.45739999999999
$ws.Range("D29").Value = -6.928600000000001
$ws.Range("D30").Value = -8.016399999999999
$ws.Range("D31").Value = -7.594099999999998
$ws.Range("E31").Value = 16.78170000000001
$ws.Range("D35").Value = -8.499299999999998
$ws.Range("A37").Value = -18.9196
$ws.Range("E38").Value = 16.24929999999999
$ws.Range("D40").Value = -8.222799999999999
$ws.Range("E42").Value = 16.3059
$ws.Range("D44").Value = -7.3085
$ws.Range("E44").Value = 16.83280000000001
$ws.Range("E47").Value = 16.35869999999998
$ws.Range("D50").Value = -7.840799999999997
$ws.Range("D54").Value = -8.229600000000003
$ws.Range("A55").Value = -22.34280000000001
$ws.Range("E56").Value = 16.4225
$ws.Range("E58").Value = 16.33370000000001
$ws.Range("E65").Value = 17.22220000000001
$ws.Range("A68").Value = -21.5448
$ws.Range("D68").Value = -6.936599999999998
$ws.Range("E73").Value = 17.21200000000001
$ws.Range("D76").Value = -7.823200000000001
$ws.Range("A77").Value = -21.35739999999998
$ws.Range("A78").Value = -21.13639999999998
$ws.Range("A79").Value = -21.70849999999999
$ws.Range("A80").Value = -20.18069999999999
$ws.Range("A81").Value = -21.76600000000001
$ws.Range("A82").Value = -22.08280000000001
$ws.Range("A84").Value = -22.08620000000001
$ws.Range("D87").Value = -8.059399999999995
$ws.Range("D88").Value = -7.143599999999998
$ws.Range("E90").Value = 16.37499999999999
$ws.Range("D92").Value = -7.048299999999999
$ws.Range("E92").Value = 17.4271
$ws.Range("E94").Value = 19.28330000000002
$ws.Range("E95").Value = 18.20770000000002
$ws.Range("D96").Value = -8.009800000000002
$ws.Range("D98").Value = -8.279400000000006
$ws.Range("A101").Value = -21.47049999999998
$ws.Range("D101").Value = -7.858100000000003
$ws.Range("E101").Value = 16.5814
$ws.Range("A102").Value = -19.1142
$ws.Range("D102").Value = -8.235599999999998
